$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.936.00"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.988.48"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.679"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.10%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.747"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000318"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "4.608.22"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "3.960.13"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "72.461.32"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  +22.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.83%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "13.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "682.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.439"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").Value = "0.0₃0857"
$ws.Range("E38").Value = "  +3.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  -3.94%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.13%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0486"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.148"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.87%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.33%  "
